$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.545.11"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "3.012.58"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.37"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.04"
$ws.Range("E6").Value = "  -2.31%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.528"
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("D9").Value = "3.010.37"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.149"
$ws.Range("E10").Value = "  -2.47%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.75"
$ws.Range("E14").Value = "  -4.15%  "
$ws.Range("E15").Value = "  +2.53%  "
$ws.Range("D16").Value = "3.514.69"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "62.537.09"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").Value = "3.013.36"
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "460.65"
$ws.Range("E20").Value = "  -3.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.04"
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.84"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  -7.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.32"
$ws.Range("E26").Value = "  -2.75%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.98"
$ws.Range("E28").Value = "  -6.20%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.02"
$ws.Range("E31").Value = "  -4.62%  "
$ws.Range("E32").Value = "  -3.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.12"
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("D35").Value = "0.0₃0815"
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("E36").Value = "  -2.50%  "
$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.12"
$ws.Range("E38").Value = "  -3.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.38"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.18"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("E41").Value = "  -10.17%  "
$ws.Range("E42").Value = "  +8.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "393.01"
$ws.Range("E43").Value = "  -9.10%  "
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("E45").Value = "  -5.71%  "
$ws.Range("D46").Value = "2.746.66"
$ws.Range("E46").Value = "  -2.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "37.57"
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.33"
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("E51").Value = "  -0.29%  "

Write-Host "Applied cryptos update"
